# Insert a new, centered contact-info paragraph right after the
# "Dheeraj Chand" name paragraph at the top of the resume.
#
# Strategy: inserting a paragraph directly after paragraph 1 (the name)
# inherits its bold/28pt run formatting, which we don't want - the new
# paragraph's run should carry no explicit run formatting at all.
# Instead we insert a brand-new blank paragraph ahead of a plain,
# unformatted paragraph (the "PROFESSIONAL SUMMARY" body text, paragraph
# 3) so it comes out clean, fill in the contact-info text and center
# alignment there, then cut it and paste it into its correct place
# right after the name paragraph.

$d = $word.ActiveDocument

$contactText = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

# Paragraph 3 is the plain "Product-focused..." summary paragraph with
# no paragraph style and no explicit run formatting - inserting before
# it yields a clean new paragraph with an empty <w:r/>.
$summaryBody = $d.Paragraphs(3)
$summaryBody.Range.InsertParagraphBefore()

$d = $word.ActiveDocument
$newPara = $d.Paragraphs(3)
$newPara.Range.Text = $contactText
$newPara.Range.ParagraphFormat.Alignment = 1

# Move the freshly-created, cleanly-formatted paragraph (with its
# paragraph mark) so it sits right after the name paragraph.
$d = $word.ActiveDocument
$moveRange = $d.Paragraphs(3).Range
$moveRange.Cut()

$d = $word.ActiveDocument
$namePara = $d.Paragraphs(1)
$insertPoint = $namePara.Range.End
$insertRange = $d.Range($insertPoint, $insertPoint)
$insertRange.Paste()
